# Coinranking crypto-price snapshot refresh (symbol list update, GitHub Actions run)
# Numeric-looking Price/Volume cells are written with a leading apostrophe so that
# Excel keeps them as plain text (matching the workbook's existing inlineStr layout)
# instead of re-interpreting "300.73" or "-4.08%" as numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.73"
$ws.Range("E2").Value = "'-4.08%"
$ws.Range("D3").Value = "'35.25"
$ws.Range("E3").Value = "'-0.67%"
$ws.Range("D4").Value = "'5.051"
$ws.Range("E4").Value = "'-0.69%"
$ws.Range("D5").Value = "'0.07963"
$ws.Range("E5").Value = "'-2.55%"
$ws.Range("E6").Value = "'-9.68%"
$ws.Range("D7").Value = "'7.777"
$ws.Range("E7").Value = "'-2.50%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.035"
$ws.Range("E8").Value = "'-2.50%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9271"
$ws.Range("E9").Value = "'-0.20%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1351"
$ws.Range("E10").Value = "'30.41%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1896"
$ws.Range("E11").Value = "'-1.32%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09045"
$ws.Range("E12").Value = "'-1.62%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03427"
$ws.Range("E13").Value = "'-5.89%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09913"
$ws.Range("E14").Value = "'0.22%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001384"
$ws.Range("E15").Value = "'-4.39%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005906"
$ws.Range("E16").Value = "'3.67%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.528"
$ws.Range("E17").Value = "'1.54%"
$ws.Range("E18").Value = "'-1.25%"
$ws.Range("E19").Value = "'-0.13%"
$ws.Range("D20").Value = "'0.1293"
$ws.Range("E20").Value = "'-0.63%"
$ws.Range("D21").Value = "'5.047"
$ws.Range("E21").Value = "'-1.11%"
$ws.Range("D22").Value = "'0.2399"
$ws.Range("E22").Value = "'8.40%"
$ws.Range("D23").Value = "'0.04496"
$ws.Range("E23").Value = "'-1.13%"
$ws.Range("D24").Value = "'0.001213"
$ws.Range("E24").Value = "'-1.56%"
$ws.Range("D25").Value = "'0.004768"
$ws.Range("E25").Value = "'-0.38%"
$ws.Range("D26").Value = "'0.0001229"
$ws.Range("E26").Value = "'-1.84%"
$ws.Range("D27").Value = "'0.0002999"
$ws.Range("E27").Value = "'-32.67%"
$ws.Range("D39").Value = "'0.01879"
$ws.Range("E39").Value = "'-5.70%"
$ws.Range("D40").Value = "'0.04762"
$ws.Range("E40").Value = "'-3.02%"
$ws.Range("D41").Value = "'0.01031"
$ws.Range("E41").Value = "'28.24%"
$ws.Range("D42").Value = "'0.007344"
$ws.Range("E42").Value = "'-2.61%"
$ws.Range("D43").Value = "'0.1321"
$ws.Range("E43").Value = "'-4.47%"
$ws.Range("D44").Value = "'0.002108"
$ws.Range("E44").Value = "'-5.19%"
$ws.Range("D45").Value = "'0.01100"
$ws.Range("E45").Value = "'-4.99%"
$ws.Range("D46").Value = "'0.00006287"
$ws.Range("E46").Value = "'-4.78%"
$ws.Range("E47").Value = "'-0.16%"
$ws.Range("D48").Value = "'64.65"
$ws.Range("E48").Value = "'9.90%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.16%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.16%"
